# Add the "ultimos tres en body" row (row 11) to the Hoja1 worksheet,
# matching the new commit: Gaizka / new task description / 2025-05-01 (45778)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("B11").Value = "Gaizka"
$ws.Range("C11").Value = "Ultimos tres en body, estilos de crear "
$ws.Range("D11").Value = 45778
$ws.Range("D11").NumberFormat = $ws.Range("D10").NumberFormat

# Update the selection to reflect where the cursor ended up after the edit
$ws.Range("D12").Select()
